$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.876.43"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "1.831.06"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.94"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6885"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07679"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3050"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07815"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.831.26"
$ws.Range("E12").Value = "  -4.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.083"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.31"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6818"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.451"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008291"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.00%  "
$ws.Range("D18").Value = "28.881.49"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.78"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.04%  "
$ws.Range("D20").Value = "2.076.31"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.70"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.464"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1479"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.20"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.795"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.19"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.546"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.215"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.147"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.180"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05100"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7656"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.838"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01849"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").Value = "1.222.01"
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.696"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9449"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.36"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.740"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000122"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "1.976.52"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.509"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.11"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.745"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4193"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.35%  "
